$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build a 22x7 array holding the new values for rows 2-23, columns A-G
$data = New-Object 'object[,]' 22,7

$data[0,0] = 5489.077381324905
$data[0,1] = 164084.9901599454
$data[0,2] = 29.893
$data[0,3] = 0
$data[0,4] = 0
$data[0,5] = 0
$data[0,6] = 0
$data[1,0] = 837.2059528526869
$data[1,1] = 26552.82400067582
$data[1,2] = 31.716
$data[1,3] = -1046.507441065859
$data[1,4] = -1.25
$data[1,5] = 0
$data[1,6] = 0
$data[2,0] = 589.7202104509649
$data[2,1] = 18703.56619466281
$data[2,2] = 31.716
$data[2,3] = -737.1502630637062
$data[2,4] = -1.25
$data[2,5] = 0
$data[2,6] = 0
$data[3,0] = 934.4551556770092
$data[3,1] = 27933.66796865284
$data[3,2] = 29.893
$data[3,3] = -1892.271690245944
$data[3,4] = -2.025
$data[3,5] = 0
$data[3,6] = 0
$data[4,0] = 1078
$data[4,1] = 30610.888
$data[4,2] = 28.396
$data[4,3] = -1347.5
$data[4,4] = -1.25
$data[4,5] = 0
$data[4,6] = 0
$data[5,0] = 18131.67667630768
$data[5,1] = 627356.0130002459
$data[5,2] = 34.6
$data[5,3] = 44603.92462371691
$data[5,4] = 2.46
$data[5,5] = 0
$data[5,6] = 0
$data[6,0] = 640.6634487247856
$data[6,1] = 38079.11340185508
$data[6,2] = 59.437
$data[6,3] = 3767.101078501739
$data[6,4] = 5.88
$data[6,5] = 0
$data[6,6] = 0
$data[7,0] = 34355.11037416998
$data[7,1] = 969329.4392072061
$data[7,2] = 28.215
$data[7,3] = 3435.511037416999
$data[7,4] = 0.1
$data[7,5] = 0
$data[7,6] = 0
$data[8,0] = 609.7648057024572
$data[8,1] = 3048.824028512286
$data[8,2] = 5
$data[8,3] = -2439.059222809829
$data[8,4] = -4
$data[8,5] = 0
$data[8,6] = 0
$data[9,0] = 3632.315283511329
$data[9,1] = 133542.071398294
$data[9,2] = 36.765
$data[9,3] = -14529.26113404532
$data[9,4] = -4
$data[9,5] = 0
$data[9,6] = 0
$data[10,0] = 874.1792747922477
$data[10,1] = 33552.74892507605
$data[10,2] = 38.382
$data[10,3] = 2149.606836714137
$data[10,4] = 2.459
$data[10,5] = 0
$data[10,6] = 0
$data[11,0] = 968.3975456496001
$data[11,1] = 31614.30627527685
$data[11,2] = 32.646
$data[11,3] = 2381.289564752367
$data[11,4] = 2.459
$data[11,5] = 0
$data[11,6] = 0
$data[12,0] = 968.3975456496001
$data[12,1] = 34689.93688025998
$data[12,2] = 35.822
$data[12,3] = 2381.289564752367
$data[12,4] = 2.459
$data[12,5] = 0
$data[12,6] = 0
$data[13,0] = 3162.2887004064
$data[13,1] = 100073.788213061
$data[13,2] = 31.646
$data[13,3] = 7776.067914299339
$data[13,4] = 2.459
$data[13,5] = 0
$data[13,6] = 0
$data[14,0] = 3162.2887004064
$data[14,1] = 110117.2171255517
$data[14,2] = 34.822
$data[14,3] = 7776.067914299339
$data[14,4] = 2.459
$data[14,5] = 0
$data[14,6] = 0
$data[15,0] = 1452.31617179842
$data[15,1] = 10166.21320258894
$data[15,2] = 7
$data[15,3] = 1452.31617179842
$data[15,4] = 1
$data[15,5] = 0
$data[15,6] = 0
$data[16,0] = 866.3
$data[16,1] = 35403.08209999999
$data[16,2] = 40.867
$data[16,3] = 1064.6827
$data[16,4] = 1.229
$data[16,5] = 0
$data[16,6] = 0
$data[17,0] = 1455.253333333333
$data[17,1] = 36547.23221333334
$data[17,2] = 25.114
$data[17,3] = -727.6266666666667
$data[17,4] = -0.5
$data[17,5] = 0
$data[17,6] = 0
$data[18,0] = 3528.131515227642
$data[18,1] = 105466.4353846999
$data[18,2] = 29.893
$data[18,3] = 3528.131515227642
$data[18,4] = 1
$data[18,5] = 0
$data[18,6] = 0
$data[19,0] = 916.3566448633733
$data[19,1] = 19243.48954213084
$data[19,2] = 21
$data[19,3] = -1145.445806079217
$data[19,4] = -1.25
$data[19,5] = 0
$data[19,6] = 0
$data[20,0] = 1728
$data[20,1] = 65342.592
$data[20,2] = 37.814
$data[20,3] = 864
$data[20,4] = 0.5
$data[20,5] = 0
$data[20,6] = 0
$data[21,0] = 600
$data[21,1] = 8194.799999999999
$data[21,2] = 13.658
$data[21,3] = 300
$data[21,4] = 0.5
$data[21,5] = 0
$data[21,6] = 0

$ws.Range("A2:G23").Value = $data

Write-Output "Edit applied: wrote $($data.Length) cells to A2:G23"
